$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "30.561.75"
Set-TextValue "E2" "  -0.79%  "
Set-TextValue "D3" "1.882.14"
Set-TextValue "E3" "  -0.43%  "
Set-TextValue "D4" "1.000"
Set-TextValue "E4" "  -0.02%  "
Set-TextValue "D5" "246.15"
Set-TextValue "E5" "  -0.55%  "
Set-TextValue "E6" "  +0.00%  "
Set-TextValue "D7" "0.4731"
Set-TextValue "E7" "  -0.12%  "
Set-TextValue "D8" "0.2888"
Set-TextValue "E8" "  -1.72%  "
Set-TextValue "D9" "0.06525"
Set-TextValue "E9" "  -0.25%  "
Set-TextValue "D10" "22.33"
Set-TextValue "E10" "  -0.60%  "
Set-TextValue "B11" "Litecoin"
Set-TextValue "C11" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D11" "100.17"
Set-TextValue "E11" "  +3.27%  "
Set-TextValue "B12" "Polygon"
Set-TextValue "C12" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D12" "0.7614"
Set-TextValue "E12" "  +2.95%  "
Set-TextValue "D13" "0.07814"
Set-TextValue "E13" "  +0.20%  "
Set-TextValue "D14" "1.880.86"
Set-TextValue "E14" "  -0.53%  "
Set-TextValue "D15" "5.236"
Set-TextValue "E15" "  -0.34%  "
Set-TextValue "D16" "284.91"
Set-TextValue "E16" "  -1.25%  "
Set-TextValue "D17" "30.547.54"
Set-TextValue "E17" "  -0.77%  "
Set-TextValue "D18" "13.20"
Set-TextValue "E18" "  -0.81%  "
Set-TextValue "D19" "0.000007518"
Set-TextValue "E19" "  -0.35%  "
Set-TextValue "E20" "  -0.07%  "
Set-TextValue "D21" "2.129.47"
Set-TextValue "E21" "  -0.48%  "
Set-TextValue "D22" "5.346"
Set-TextValue "E22" "  +0.00%  "
Set-TextValue "D23" "1.001"
Set-TextValue "E23" "  +0.03%  "
Set-TextValue "D24" "6.435"
Set-TextValue "E24" "  +2.34%  "
Set-TextValue "D25" "9.162"
Set-TextValue "E25" "  -0.68%  "
Set-TextValue "D26" "163.27"
Set-TextValue "E26" "  -0.97%  "
Set-TextValue "D27" "19.04"
Set-TextValue "E27" "  +0.00%  "
Set-TextValue "D28" "1.903"
Set-TextValue "E28" "  -0.97%  "
Set-TextValue "D29" "0.09684"
Set-TextValue "E29" "  -0.88%  "
Set-TextValue "D30" "1.327"
Set-TextValue "E30" "  -0.99%  "
Set-TextValue "D31" "1.502"
Set-TextValue "E31" "  +0.81%  "
Set-TextValue "D32" "4.244"
Set-TextValue "E32" "  -1.98%  "
Set-TextValue "D33" "4.171"
Set-TextValue "E33" "  +0.04%  "
Set-TextValue "D34" "0.04839"
Set-TextValue "E34" "  -1.15%  "
Set-TextValue "D35" "1.127"
Set-TextValue "E35" "  -0.33%  "
Set-TextValue "D36" "0.6978"
Set-TextValue "E36" "  -0.47%  "
Set-TextValue "D37" "2.760"
Set-TextValue "E37" "  +1.14%  "
Set-TextValue "D38" "0.01903"
Set-TextValue "E38" "  -0.10%  "
Set-TextValue "D39" "2.884"
Set-TextValue "E39" "  +2.34%  "
Set-TextValue "B40" "FraxShare"
Set-TextValue "C40" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D40" "6.289"
Set-TextValue "E40" "  -0.80%  "
Set-TextValue "B41" "Aave"
Set-TextValue "C41" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D41" "75.49"
Set-TextValue "E41" "  -1.53%  "
Set-TextValue "D42" "1.975"
Set-TextValue "E42" "  -1.32%  "
Set-TextValue "D43" "0.4245"
Set-TextValue "E43" "  -1.38%  "
Set-TextValue "E44" "  -0.05%  "
Set-TextValue "D45" "0.8386"
Set-TextValue "E45" "  +0.06%  "
Set-TextValue "D46" "101.21"
Set-TextValue "E46" "  -0.82%  "
Set-TextValue "D47" "9.866"
Set-TextValue "E47" "  +3.09%  "
Set-TextValue "D48" "7.016"
Set-TextValue "E48" "  -0.60%  "
Set-TextValue "D49" "35.22"
Set-TextValue "E49" "  -1.11%  "
Set-TextValue "D50" "0.05782"
Set-TextValue "E50" "  +0.39%  "
Set-TextValue "D51" "0.3954"
Set-TextValue "E51" "  -0.76%  "
